$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: date update (leading apostrophe keeps it stored as text, matching
# the original cell which holds the literal string "01/01/2018", not a date)
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"

# Programa resumido (PT)
$ws.Range("B14").Value = "Propriedades físicas, acidez e basicidade dos compostos orgânicos. Classes de compostos orgânicos, estereoquimica, principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Importância Industrial, impacto ambiental e degradação de compostos orgânicos."
$ws.Range("C14").Value = "Propriedades físicas, acidez e basicidade dos compostos orgânicos. Classes de compostos orgânicos, estereoquimica, principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Importância Industrial, impacto ambiental e degradação de compostos orgânicos."

# Short syllabus (EN)
$ws.Range("B15").Value = "Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, stereochemistry, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds."
$ws.Range("C15").Value = "Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, stereochemistry, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds."

# Programa (PT)
$ws.Range("B16").Value = "Terminologia e conceitos básicos de química orgânica, estrutura e características das importantes classes dos compostos orgânicos, relacionando com produtos orgânicos antropogênicos. Acidez e basicidade em compostos orgânicos. Discussão dos aspectos básicos dos processos de partição. Como a estrutura química determina a solubilidade do composto orgânico em água. Estereoquimica. Principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Estereoquímica. Como os absorventes naturais são importantes para transporte, distribuição e destino dos compostos orgânicos no meio ambiente. Processos de transformação abiótica e biológicas. Processos químicos, fotoquímicos, biológicos e reações de degradação. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aquíferos."
$ws.Range("C16").Value = "Terminologia e conceitos básicos de química orgânica, estrutura e características das importantes classes dos compostos orgânicos, relacionando com produtos orgânicos antropogênicos. Acidez e basicidade em compostos orgânicos. Discussão dos aspectos básicos dos processos de partição. Como a estrutura química determina a solubilidade do composto orgânico em água. Estereoquimica. Principais reações e introdução à mecanismos de reações orgânicas (Substituição, eliminação, adição; via iônica e radicalar). Estereoquímica. Como os absorventes naturais são importantes para transporte, distribuição e destino dos compostos orgânicos no meio ambiente. Processos de transformação abiótica e biológicas. Processos químicos, fotoquímicos, biológicos e reações de degradação. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aquíferos."

# Syllabus (EN)
$ws.Range("B17").Value = "Terminology and basic concepts of organic chemistry, structure and characteristics of the important classes of organic compounds, relating to anthropogenic organic products. Acidity and basicity in organic compounds. Discussion of the basic aspects of partition processes. How the chemical structure determines the solubility of the organic compound in water. Stereochemistry. Main reactions and introduction to the organic reaction mechanisms (Substitution, elimination, addition; ionic and radical pathways). Stereochemistry. As natural absorbents are important for the transport, distribution and destination of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the discussed concepts and environmental systems, such as lakes, rivers and aquifers."
$ws.Range("C17").Value = "Terminology and basic concepts of organic chemistry, structure and characteristics of the important classes of organic compounds, relating to anthropogenic organic products. Acidity and basicity in organic compounds. Discussion of the basic aspects of partition processes. How the chemical structure determines the solubility of the organic compound in water. Stereochemistry. Main reactions and introduction to the organic reaction mechanisms (Substitution, elimination, addition; ionic and radical pathways). Stereochemistry. As natural absorbents are important for the transport, distribution and destination of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the discussed concepts and environmental systems, such as lakes, rivers and aquifers."

# Bibliografia
$ws.Range("B22").Value = "- Bruice, Paula Yurkanis – Química Orgânica – PEARSON Prentice Hall – São Paulo, 2006. Solomons, T.W.G.- Química Orgânica. Volumes 1 e 2, Rio de Janeiro, Livros Técnicos e Científicos, 2012.- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry – Wiley Interscience, 2016.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010."
$ws.Range("C22").Value = "- Bruice, Paula Yurkanis – Química Orgânica – PEARSON Prentice Hall – São Paulo, 2006. Solomons, T.W.G.- Química Orgânica. Volumes 1 e 2, Rio de Janeiro, Livros Técnicos e Científicos, 2012.- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry – Wiley Interscience, 2016.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010."
